$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = "31/07/2021"
$ws.Range("C32").Value = "QMVAR"
$ws.Range("D32").Value = "Weekly Graph dynamic code generation reports checked with JP reports and mismatching datas highlighted in red and sent to Mohan san"

$ws.Range("D32").Select()
